# The diff replaces the per-row price/variety/date details (columns D and
# K through T) of the weekly "Damasco" price sheet with the values from a
# different week, row by row. Comparing old vs. new values shows this is a
# permutation of the 15 data rows (rows 2-16): each destination row's new
# content is exactly the old content of another row in the same block.
#
# destination row -> source row (of the values as they existed before edit)
#   2 <- 5    3 <- 10   4 <- 14   5 <- 13   6 <- 8    7 <- 12   8 <- 6
#   9 <- 7    10 <- 11  11 <- 15  12 <- 16  13 <- 9   14 <- 3   15 <- 4
#   16 <- 2
#
# Columns A, B, C, E, F, G, H, I, J are identical in every data row, so we
# copy the whole row (columns A-T) for simplicity and safety; this is
# equivalent to only touching D and K:T as the diff shows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstDataRow = 2
$lastDataRow = 16
$lastCol = 20   # column T

$rowMap = @{
    2  = 5
    3  = 10
    4  = 14
    5  = 13
    6  = 8
    7  = 12
    8  = 6
    9  = 7
    10 = 11
    11 = 15
    12 = 16
    13 = 9
    14 = 3
    15 = 4
    16 = 2
}

# Snapshot every source row first so that writing destination rows never
# clobbers data that is still needed as a source later on (this matters
# because the mapping above is a full permutation with no fixed points).
$snapshot = @{}
for ($r = $firstDataRow; $r -le $lastDataRow; $r++) {
    $values = @()
    for ($c = 1; $c -le $lastCol; $c++) {
        $values += , ($ws.Cells.Item($r, $c).Value2)
    }
    $snapshot[$r] = $values
}

foreach ($destRow in $rowMap.Keys | Sort-Object) {
    $srcRow = $rowMap[$destRow]
    $values = $snapshot[$srcRow]
    for ($c = 1; $c -le $lastCol; $c++) {
        $ws.Cells.Item($destRow, $c).Value = $values[$c - 1]
    }
}
